$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters for the 20 data columns (A..T)
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Updated/added rows 2-10 (rows 2-7 revised in place, rows 8-10 newly appended)
$data = @{
  2  = @("ECs","Fgf1","Cd44","ECs",  2, 0.6666666666666666, 2.004760666666666, 6.014282, 0.1200698528618338, 0.1200698528618338, 3, 1, 261.380203, 784.1406089999999, 0.6968677182772199, 0.6968677182772199, 524.0047500197486, 4716.042750177738, 0.08367280439770763, 0.08367280439770763)
  3  = @("ECs","Fgf1","Cd44","FAPs", 2, 0.6666666666666666, 2.004760666666666, 6.014282, 0.1200698528618338, 0.1200698528618338, 3, 1, 31.999428, 95.998284, 0.08531391482826334, 0.08531391482826335, 64.15119461023198, 577.360751492088, 0.0102436292004966, 0.0102436292004966)
  4  = @("ECs","Fgf1","Cd44","sCs",  2, 0.6666666666666666, 2.004760666666666, 6.014282, 0.1200698528618338, 0.1200698528618338, 3, 1, 81.699019, 245.097057, 0.2178183668945166, 0.2178183668945167, 163.7869797964527, 1474.082818168074, 0.02615341926362954, 0.02615341926362954)
  5  = @("FAPs","Fgf1","Cd44","ECs", 3, 1, 9.409654999999999, 28.228965, 0.5635664696121425, 0.5635664696121424, 3, 1, 261.380203, 784.1406089999999, 0.6968677182772199, 0.6968677182772199, 2459.497534059965, 22135.47780653968, 0.392731279776162, 0.3927312797761619)
  6  = @("FAPs","Fgf1","Cd44","FAPs",3, 1, 9.409654999999999, 28.228965, 0.5635664696121425, 0.5635664696121424, 3, 1, 31.999428, 95.998284, 0.08531391482826334, 0.08531391482826335, 301.10357767734, 2709.93219909606, 0.04808006178855539, 0.04808006178855538)
  7  = @("FAPs","Fgf1","Cd44","sCs", 3, 1, 9.409654999999999, 28.228965, 0.5635664696121425, 0.5635664696121424, 3, 1, 81.699019, 245.097057, 0.2178183668945166, 0.2178183668945167, 768.759582628445, 6918.836243656005, 0.1227551280474251, 0.1227551280474251)
  8  = @("sCs","Fgf1","Cd44","ECs",  3, 1, 5.282203999999999, 15.846612, 0.3163636775260238, 0.3163636775260238, 3, 1, 261.380203, 784.1406089999999, 0.6968677182772199, 0.6968677182772199, 1380.663553807412, 12425.97198426671, 0.2204636341033504, 0.2204636341033504)
  9  = @("sCs","Fgf1","Cd44","FAPs", 3, 1, 5.282203999999999, 15.846612, 0.3163636775260238, 0.3163636775260238, 3, 1, 31.999428, 95.998284, 0.08531391482826334, 0.08531391482826335, 169.027506579312, 1521.247559213808, 0.02699022383921136, 0.02699022383921137)
  10 = @("sCs","Fgf1","Cd44","sCs",  3, 1, 5.282203999999999, 15.846612, 0.3163636775260238, 0.3163636775260238, 3, 1, 81.699019, 245.097057, 0.2178183668945166, 0.2178183668945167, 431.550884957876, 3883.957964620884, 0.068909819583462, 0.06890981958346201)
}

foreach ($r in 2..10) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}
